$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("HVAC_ASSEMBLIES")

# Rows whose heating/cooling season start/end columns (G,H,I,J) need updating.
# (Residential "MULTIRES_*" rows that are not the AC-equipped variants.)
$rows = @(4,5,6,7,8,9,10,11,15,16,17,18,19,20,21)

foreach ($r in $rows) {
    # Set column I ("00|00") first so it becomes the first newly-introduced
    # shared string, then G ("01|01"), then H ("31|12"), matching the order
    # the strings were appended to the shared-strings table.
    $ws.Range("I$r").Value = "00|00"
    $ws.Range("G$r").Value = "01|01"
    $ws.Range("H$r").Value = "31|12"
    $ws.Range("J$r").Value = "00|00"
}

# Update the active sheet / selection state: ENVELOPE_ASSEMBLIES loses the
# tab-selected flag and its remembered selection moves to N4; HVAC_ASSEMBLIES
# becomes the active/selected tab with its selection at H23.
$wsEnvelope = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")
[void]$wsEnvelope.Activate()
[void]$wsEnvelope.Range("N4").Select()

[void]$ws.Activate()
[void]$ws.Range("H23").Select()
